# Apply "Penalty Reward System" updates to the forecast workbook.
# Sheet 1 "Forecast Comparison": shift Week_Start_Date (col B) forward by one
# week and update MyForecast (col D) values for rows 2-17.
# Sheet 2 "Summary": update several derived metric values in col B.
#
# NOTE: Week_Start_Date and all Summary "Value" cells are stored as plain
# TEXT (not real dates/numbers) in the source workbook. A leading apostrophe
# forces Excel to keep the literal text instead of auto-converting
# date-shaped / number-shaped strings into typed date or numeric cells.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Sheet 1: Forecast Comparison ---

$wsForecast.Range("B2").Value  = "'2025-01-12"
$wsForecast.Range("D2").Value  = 57

$wsForecast.Range("B3").Value  = "'2025-01-19"
$wsForecast.Range("D3").Value  = 58

$wsForecast.Range("B4").Value  = "'2025-01-26"
$wsForecast.Range("D4").Value  = 58

$wsForecast.Range("B5").Value  = "'2025-02-02"
$wsForecast.Range("D5").Value  = 59

$wsForecast.Range("B6").Value  = "'2025-02-09"
$wsForecast.Range("D6").Value  = 62

$wsForecast.Range("B7").Value  = "'2025-02-16"
$wsForecast.Range("D7").Value  = 65

$wsForecast.Range("B8").Value  = "'2025-02-23"
$wsForecast.Range("D8").Value  = 67

$wsForecast.Range("B9").Value  = "'2025-03-02"
$wsForecast.Range("D9").Value  = 66

$wsForecast.Range("B10").Value = "'2025-03-09"
$wsForecast.Range("D10").Value = 67

$wsForecast.Range("B11").Value = "'2025-03-16"

$wsForecast.Range("B12").Value = "'2025-03-23"
$wsForecast.Range("D12").Value = 48

$wsForecast.Range("B13").Value = "'2025-03-30"
$wsForecast.Range("D13").Value = 48

$wsForecast.Range("B14").Value = "'2025-04-06"
$wsForecast.Range("D14").Value = 47

$wsForecast.Range("B15").Value = "'2025-04-13"
$wsForecast.Range("D15").Value = 46

$wsForecast.Range("B16").Value = "'2025-04-20"
$wsForecast.Range("D16").Value = 46

$wsForecast.Range("B17").Value = "'2025-04-27"
$wsForecast.Range("D17").Value = 45

# --- Sheet 2: Summary ---

$wsSummary.Range("B2").Value  = "'2024-02-04 to 2025-01-05"
$wsSummary.Range("B4").Value  = "'94"
$wsSummary.Range("B6").Value  = "'44"
$wsSummary.Range("B8").Value  = "'1862 units"
$wsSummary.Range("B9").Value  = "'889"
$wsSummary.Range("B10").Value = "'492"
$wsSummary.Range("B11").Value = "'232"
$wsSummary.Range("B12").Value = "'67"
$wsSummary.Range("B13").Value = "'2025-02-23"
$wsSummary.Range("B15").Value = "'2025-04-27"
